$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new numeric-looking text must be forced to Text format
# so Excel does not auto-convert them to numbers (losing the exact
# decimal-formatted representation used throughout the sheet).
$textCells = @("D5", "D6", "D7", "D10", "D12", "D13", "D14", "D16", "D19", "D22", "D24", "D28", "D29", "D30", "D32", "D33", "D34", "D35", "D36", "D38", "D39", "D40", "D46", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values scraped by the GitHub Actions job
$ws.Cells.Item(2, 4).Value = '36.961.70'
$ws.Cells.Item(2, 5).Value = '  +0.84%  '
$ws.Cells.Item(3, 4).Value = '1.985.01'
$ws.Cells.Item(3, 5).Value = '  +1.05%  '
$ws.Cells.Item(4, 5).Value = '  -0.05%  '
$ws.Cells.Item(5, 4).Value = '245.30'
$ws.Cells.Item(5, 5).Value = '  +0.33%  '
$ws.Cells.Item(6, 4).Value = '0.629'
$ws.Cells.Item(6, 5).Value = '  +1.62%  '
$ws.Cells.Item(7, 4).Value = '61.03'
$ws.Cells.Item(7, 5).Value = '  +3.00%  '
$ws.Cells.Item(8, 5).Value = '  -0.01%  '
$ws.Cells.Item(9, 5).Value = '  +2.10%  '
$ws.Cells.Item(10, 4).Value = '0.0801'
$ws.Cells.Item(10, 5).Value = '  -1.56%  '
$ws.Cells.Item(11, 5).Value = '  +0.73%  '
$ws.Cells.Item(12, 4).Value = '14.96'
$ws.Cells.Item(12, 5).Value = '  +9.02%  '
$ws.Cells.Item(13, 2).Value = 'Avalanche'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Cells.Item(13, 4).Value = '22.19'
$ws.Cells.Item(13, 5).Value = '  -0.48%  '
$ws.Cells.Item(14, 2).Value = 'Polygon'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Cells.Item(14, 4).Value = '0.846'
$ws.Cells.Item(14, 5).Value = '  +2.03%  '
$ws.Cells.Item(15, 4).Value = '2.271.09'
$ws.Cells.Item(15, 5).Value = '  +0.80%  '
$ws.Cells.Item(16, 4).Value = '5.48'
$ws.Cells.Item(16, 5).Value = '  +4.20%  '
$ws.Cells.Item(17, 4).Value = '1.979.93'
$ws.Cells.Item(17, 5).Value = '  +0.65%  '
$ws.Cells.Item(18, 4).Value = '36.839.75'
$ws.Cells.Item(18, 5).Value = '  +0.72%  '
$ws.Cells.Item(19, 4).Value = '70.29'
$ws.Cells.Item(19, 5).Value = '  +0.40%  '
$ws.Cells.Item(20, 5).Value = '  +0.47%  '
$ws.Cells.Item(21, 5).Value = '  +2.39%  '
$ws.Cells.Item(22, 4).Value = '229.87'
$ws.Cells.Item(22, 5).Value = '  +0.28%  '
$ws.Cells.Item(23, 5).Value = '  +0.04%  '
$ws.Cells.Item(24, 4).Value = '2.51'
$ws.Cells.Item(24, 5).Value = '  +2.58%  '
$ws.Cells.Item(25, 5).Value = '  +0.86%  '
$ws.Cells.Item(26, 5).Value = '  +3.69%  '
$ws.Cells.Item(27, 5).Value = '  +0.92%  '
$ws.Cells.Item(28, 4).Value = '163.20'
$ws.Cells.Item(28, 5).Value = '  +1.92%  '
$ws.Cells.Item(29, 4).Value = '19.56'
$ws.Cells.Item(29, 5).Value = '  +0.90%  '
$ws.Cells.Item(30, 4).Value = '1.38'
$ws.Cells.Item(30, 5).Value = '  +18.84%  '
$ws.Cells.Item(31, 5).Value = '  +1.79%  '
$ws.Cells.Item(32, 4).Value = '4.90'
$ws.Cells.Item(32, 5).Value = '  +3.92%  '
$ws.Cells.Item(33, 4).Value = '0.0622'
$ws.Cells.Item(33, 5).Value = '  +0.54%  '
$ws.Cells.Item(34, 4).Value = '4.55'
$ws.Cells.Item(34, 5).Value = '  +6.29%  '
$ws.Cells.Item(35, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Cells.Item(35, 4).Value = '2.28'
$ws.Cells.Item(35, 5).Value = '  +1.31%  '
$ws.Cells.Item(36, 2).Value = 'BinanceUSD'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Cells.Item(36, 4).Value = '1.00'
$ws.Cells.Item(36, 5).Value = '  -0.08%  '
$ws.Cells.Item(37, 5).Value = '  -1.27%  '
$ws.Cells.Item(38, 4).Value = '1.78'
$ws.Cells.Item(38, 5).Value = '  -0.15%  '
$ws.Cells.Item(39, 4).Value = '5.57'
$ws.Cells.Item(39, 5).Value = '  -7.21%  '
$ws.Cells.Item(40, 4).Value = '0.0997'
$ws.Cells.Item(40, 5).Value = '  +0.84%  '
$ws.Cells.Item(41, 5).Value = '  +0.72%  '
$ws.Cells.Item(42, 5).Value = '  +0.94%  '
$ws.Cells.Item(43, 5).Value = '  +0.71%  '
$ws.Cells.Item(44, 5).Value = '  +2.83%  '
$ws.Cells.Item(45, 4).Value = '1.372.84'
$ws.Cells.Item(45, 5).Value = '  +0.97%  '
$ws.Cells.Item(46, 4).Value = '90.28'
$ws.Cells.Item(46, 5).Value = '  +2.78%  '
$ws.Cells.Item(47, 5).Value = '  +0.05%  '
$ws.Cells.Item(48, 4).Value = '7.26'
$ws.Cells.Item(48, 5).Value = '  +1.81%  '
$ws.Cells.Item(49, 4).Value = '2.82'
$ws.Cells.Item(49, 5).Value = '  -0.54%  '
$ws.Cells.Item(50, 4).Value = '46.34'
$ws.Cells.Item(50, 5).Value = '  +5.75%  '
$ws.Cells.Item(51, 4).Value = '1.99'
$ws.Cells.Item(51, 5).Value = '  +11.81%  '
